# Apply the weekly fruit/vegetable price update (Higo - Vega Central Mapocho de Santiago).
# The diff effectively rotates the date/volume/price/origin data among the
# "Primera"/"Segunda" row-pairs 2-3, 4-5, 6-7, 10-11, 14-15 while the
# descriptive columns (market, product, variety, quality, unit, etc.) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Primera)
$ws.Cells.Item(2, 4).Value  = 44320
$ws.Cells.Item(2, 13).Value = 20
$ws.Cells.Item(2, 14).Value = 12000
$ws.Cells.Item(2, 15).Value = 12000
$ws.Cells.Item(2, 16).Value = 12000
$ws.Cells.Item(2, 18).Value = "Región Metropolitana"
$ws.Cells.Item(2, 19).Value = 1714

# Row 3 (Segunda)
$ws.Cells.Item(3, 4).Value  = 44320
$ws.Cells.Item(3, 13).Value = 30
$ws.Cells.Item(3, 14).Value = 8000
$ws.Cells.Item(3, 15).Value = 8000
$ws.Cells.Item(3, 16).Value = 8000
$ws.Cells.Item(3, 18).Value = "Región Metropolitana"
$ws.Cells.Item(3, 19).Value = 1143

# Row 4 (Primera)
$ws.Cells.Item(4, 4).Value  = 44302
$ws.Cells.Item(4, 13).Value = 50
$ws.Cells.Item(4, 14).Value = 15000
$ws.Cells.Item(4, 15).Value = 15000
$ws.Cells.Item(4, 16).Value = 15000
$ws.Cells.Item(4, 19).Value = 2143

# Row 5 (Segunda)
$ws.Cells.Item(5, 4).Value  = 44302
$ws.Cells.Item(5, 14).Value = 12000
$ws.Cells.Item(5, 15).Value = 12000
$ws.Cells.Item(5, 16).Value = 12000
$ws.Cells.Item(5, 19).Value = 1714

# Row 6 (Primera)
$ws.Cells.Item(6, 4).Value  = 44292
$ws.Cells.Item(6, 13).Value = 25
$ws.Cells.Item(6, 14).Value = 16000
$ws.Cells.Item(6, 15).Value = 16000
$ws.Cells.Item(6, 16).Value = 16000
$ws.Cells.Item(6, 19).Value = 2286

# Row 7 (Segunda)
$ws.Cells.Item(7, 4).Value  = 44292
$ws.Cells.Item(7, 13).Value = 30
$ws.Cells.Item(7, 14).Value = 15000
$ws.Cells.Item(7, 15).Value = 15000
$ws.Cells.Item(7, 16).Value = 15000
$ws.Cells.Item(7, 19).Value = 2143

# Row 10 (Primera)
$ws.Cells.Item(10, 4).Value  = 44299
$ws.Cells.Item(10, 13).Value = 80
$ws.Cells.Item(10, 14).Value = 15000
$ws.Cells.Item(10, 15).Value = 15000
$ws.Cells.Item(10, 16).Value = 15000
$ws.Cells.Item(10, 18).Value = "Provincia de Santiago"
$ws.Cells.Item(10, 19).Value = 2143

# Row 11 (Segunda)
$ws.Cells.Item(11, 4).Value  = 44299
$ws.Cells.Item(11, 13).Value = 75
$ws.Cells.Item(11, 14).Value = 12000
$ws.Cells.Item(11, 15).Value = 12000
$ws.Cells.Item(11, 16).Value = 12000
$ws.Cells.Item(11, 18).Value = "Provincia de Santiago"
$ws.Cells.Item(11, 19).Value = 1714

# Row 14 (Primera)
$ws.Cells.Item(14, 4).Value  = 44322
$ws.Cells.Item(14, 13).Value = 45
$ws.Cells.Item(14, 14).Value = 12000
$ws.Cells.Item(14, 15).Value = 12000
$ws.Cells.Item(14, 16).Value = 12000
$ws.Cells.Item(14, 19).Value = 1714

# Row 15 (Segunda)
$ws.Cells.Item(15, 4).Value  = 44322
$ws.Cells.Item(15, 13).Value = 80
$ws.Cells.Item(15, 14).Value = 8000
$ws.Cells.Item(15, 15).Value = 8000
$ws.Cells.Item(15, 16).Value = 8000
$ws.Cells.Item(15, 19).Value = 1143
